$wb = $excel.ActiveWorkbook

# --- Sheet "TAC" ---
$ws = $wb.Worksheets.Item("TAC")

# Row 6: convert B6/C6 from text to real numbers
$ws.Cells.Item(6, 2).Value = 33225726842
$ws.Cells.Item(6, 3).Value = 11995037548

# Row 7: Renan rosa
$ws.Cells.Item(7, 1).Value = "Renan rosa"
$ws.Cells.Item(7, 2).Value = 48267379800
$ws.Cells.Item(7, 3).Value = 11978813206
$ws.Cells.Item(7, 4).Value = "Atibaia-SP"
$ws.Cells.Item(7, 5).Value = "TAC"
$ws.Cells.Item(7, 6).Value = "Não"
$ws.Cells.Item(7, 8).Value = "23/04/2025 13:56"
$ws.Cells.Item(7, 9).Value = "Completo"

# Row 8: Fabiano felski
$ws.Cells.Item(8, 1).Value = "Fabiano felski"
$ws.Cells.Item(8, 2).Value = 1238946976
$ws.Cells.Item(8, 3).Value = 49991863317
$ws.Cells.Item(8, 4).Value = "são miguel do oeste-sc"
$ws.Cells.Item(8, 5).Value = "TAC"
$ws.Cells.Item(8, 6).Value = "Não"
$ws.Cells.Item(8, 8).Value = "23/04/2025 14:13"
$ws.Cells.Item(8, 9).Value = "Completo"

# Row 9: Gustavo Tedesco Bedin (CPF/Telefone kept as TEXT, not number)
$ws.Cells.Item(9, 1).Value = "Gustavo Tedesco Bedin"
$ws.Cells.Item(9, 2).Value = "'10854057960"
$ws.Cells.Item(9, 3).Value = "'49991878706"
$ws.Cells.Item(9, 4).Value = "São Miguel Do Oeste-SC"
$ws.Cells.Item(9, 5).Value = "TAC"
$ws.Cells.Item(9, 6).Value = "Sim"
$ws.Cells.Item(9, 8).Value = "23/04/2025 15:11"
$ws.Cells.Item(9, 9).Value = "Completo"

# --- Sheet "Contatos Incompletos" ---
$ws2 = $wb.Worksheets.Item("Contatos Incompletos")

# Row 4: convert B4/C4 from text to real numbers
$ws2.Cells.Item(4, 2).Value = 35263195885
$ws2.Cells.Item(4, 3).Value = 19994564565

# Row 5: Thaina Cristina Ramos de Oliveira
$ws2.Cells.Item(5, 1).Value = "Thaina Cristina Ramos de Oliveira"
$ws2.Cells.Item(5, 2).Value = 45693627888
$ws2.Cells.Item(5, 3).Value = 19995305329
$ws2.Cells.Item(5, 4).Value = "Mogi Mirim-SP"
$ws2.Cells.Item(5, 5).Value = "AGREGADO"
$ws2.Cells.Item(5, 6).Value = "FTW-7533"
$ws2.Cells.Item(5, 7).Value = "'"
$ws2.Cells.Item(5, 8).Value = "23/04/2025 14:48"
$ws2.Cells.Item(5, 9).Value = "Em andamento"

# Row 6: Rodrigo Rocha de Castro (CPF/Telefone kept as TEXT, not number)
$ws2.Cells.Item(6, 1).Value = "Rodrigo Rocha de Castro"
$ws2.Cells.Item(6, 2).Value = "'10003477967"
$ws2.Cells.Item(6, 3).Value = "'49988681357"
$ws2.Cells.Item(6, 4).Value = "Caçador-sc"
$ws2.Cells.Item(6, 5).Value = "AGREGADO"
$ws2.Cells.Item(6, 6).Value = "QJI-9564"
$ws2.Cells.Item(6, 7).Value = "'"
$ws2.Cells.Item(6, 8).Value = "23/04/2025 15:17"
$ws2.Cells.Item(6, 9).Value = "Em andamento"
